$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data keyed by row number, rebuilt from the permuted/updated weekly records
# Columns: D=Fecha, K=Variedad, L=Calidad, M=Volumen, N=Precio minimo, O=Precio maximo,
#          P=Precio promedio ponderado, Q=Unidad comercializacion, R=Origen, S=Precio $/Kg, T=Kg/unidad
$rows = @(
    @{ Row=2; D=44181; K='Modesto'; L='Primera'; M=50; N=20000; O=21000; P=20500; Q='$/caja 18 kilos'; R='Región de Coquimbo'; S=1139; T=18 }
    @{ Row=3; D=44176; K='Castle Brite'; L='Primera'; M=50; N=17000; O=18000; P=17400; Q='$/caja 18 kilos'; R='Región de O''Higgins'; S=967; T=18 }
    @{ Row=4; D=44544; K='Castle Brite'; L='Segunda'; M=160; N=16000; O=17000; P=16500; Q='$/caja 15 kilos'; R='Región de O''Higgins'; S=1100; T=15 }
    @{ Row=5; D=44168; K='Castle Brite'; L='Primera'; M=30; N=16000; O=17000; P=16500; Q='$/caja 16 kilos granel'; R='Región de Coquimbo'; S=1031; T=16 }
    @{ Row=6; D=44904; K='Castle Brite'; L='Primera'; M=60; N=15000; O=16000; P=15500; Q='$/bandeja 10 kilos'; R='Región de O''Higgins'; S=1550; T=10 }
    @{ Row=7; D=44904; K='Castle Brite'; L='Segunda'; M=30; N=14000; O=14000; P=14000; Q='$/bandeja 10 kilos'; R='Región de O''Higgins'; S=1400; T=10 }
    @{ Row=8; D=44174; K='Castle Brite'; L='Primera'; M=75; N=9000; O=10000; P=9467; Q='$/caja 10 kilos'; R='Región de O''Higgins'; S=947; T=10 }
    @{ Row=9; D=44551; K='Castle Brite'; L='Primera'; M=120; N=15500; O=16000; P=15750; Q='$/caja 15 kilos'; R='Región de O''Higgins'; S=1050; T=15 }
    @{ Row=10; D=44552; K='Castle Brite'; L='Primera'; M=120; N=15500; O=16000; P=15750; Q='$/caja 15 kilos'; R='Región de O''Higgins'; S=1050; T=15 }
    @{ Row=11; D=44165; K='Castle Brite'; L='Segunda'; M=60; N=16000; O=17000; P=16500; Q='$/caja 15 kilos granel'; R='Provincia de Limarí'; S=1100; T=15 }
    @{ Row=12; D=44187; K='Dina'; L='Primera'; M=55; N=15000; O=16000; P=15455; Q='$/caja 15 kilos granel'; R='Región de O''Higgins'; S=1030; T=15 }
    @{ Row=13; D=44537; K='Castle Brite'; L='Primera'; M=60; N=21000; O=21500; P=21250; Q='$/caja 15 kilos'; R='Región de O''Higgins'; S=1417; T=15 }
    @{ Row=14; D=44189; K='Dina'; L='Primera'; M=80; N=16000; O=17000; P=16562; Q='$/caja 18 kilos'; R='Región de O''Higgins'; S=920; T=18 }
    @{ Row=15; D=44907; K='Castle Brite'; L='Primera'; M=120; N=15000; O=16000; P=15500; Q='$/bandeja 10 kilos'; R='Región de O''Higgins'; S=1550; T=10 }
    @{ Row=16; D=44907; K='Castle Brite'; L='Segunda'; M=60; N=14000; O=14000; P=14000; Q='$/bandeja 10 kilos'; R='Región de O''Higgins'; S=1400; T=10 }
)

foreach ($rec in $rows) {
    $r = $rec.Row
    $ws.Range("D$r").Value = $rec.D
    $ws.Range("K$r").Value = $rec.K
    $ws.Range("L$r").Value = $rec.L
    $ws.Range("M$r").Value = $rec.M
    $ws.Range("N$r").Value = $rec.N
    $ws.Range("O$r").Value = $rec.O
    $ws.Range("P$r").Value = $rec.P
    $ws.Range("Q$r").Value = $rec.Q
    $ws.Range("R$r").Value = $rec.R
    $ws.Range("S$r").Value = $rec.S
    $ws.Range("T$r").Value = $rec.T
}
